$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 14
$ws.Range("B2").Value = "KEGG_2021_Human: PI3K-Akt signaling pathway"
$ws.Range("C2").Value = [double]"5.558218211829239e-12"
$ws.Range("D2").Value = [double]"4.921079877004785e-08"
$ws.Range("E2").Value = [double]"6.12430386631518e-07"
$ws.Range("A3").Value = 20
$ws.Range("B3").Value = "KEGG_2021_Human: Focal adhesion"
$ws.Range("C3").Value = [double]"1.984119224391844e-05"
$ws.Range("D3").Value = [double]"2.697146656837494e-05"
$ws.Range("E3").Value = [double]"4.272537005773376e-06"
$ws.Range("A4").Value = 37
$ws.Range("B4").Value = "KEGG_2021_Human: TGF-beta signaling pathway"
$ws.Range("C4").Value = [double]"2.140821594454763e-06"
$ws.Range("D4").Value = [double]"0.002296359605154307"
$ws.Range("E4").Value = [double]"0.002229718270090995"
$ws.Range("A5").Value = 42
$ws.Range("B5").Value = "GO_Cellular_Component_2018: focal adhesion (GO:0005925)"
$ws.Range("C5").Value = [double]"9.247551505658709e-07"
$ws.Range("D5").Value = [double]"1.845273399229044e-11"
$ws.Range("E5").Value = "ns"
$ws.Range("A6").Value = 117
$ws.Range("B6").Value = "KEGG_2021_Human: Tight junction"
$ws.Range("C6").Value = [double]"0.0003439937129393427"
$ws.Range("D6").Value = [double]"6.580526896332986e-05"
$ws.Range("E6").Value = "ns"
$ws.Range("A7").Value = 154
$ws.Range("B7").Value = "KEGG_2021_Human: Adherens junction"
$ws.Range("C7").Value = [double]"6.077939248200902e-07"
$ws.Range("D7").Value = [double]"0.0001202543103596818"
$ws.Range("E7").Value = "ns"
$ws.Range("A8").Value = 160
$ws.Range("B8").Value = "KEGG_2021_Human: Regulation of actin cytoskeleton"
$ws.Range("C8").Value = [double]"7.658662304895348e-07"
$ws.Range("D8").Value = [double]"8.427888745314609e-07"
$ws.Range("E8").Value = "ns"
$ws.Range("A9").Value = 164
$ws.Range("B9").Value = "KEGG_2021_Human: Hippo signaling pathway"
$ws.Range("C9").Value = [double]"0.002141708366094922"
$ws.Range("D9").Value = [double]"0.0004380982536883778"
$ws.Range("E9").Value = "ns"
$ws.Range("A10").Value = 242
$ws.Range("B10").Value = "GO_Biological_Process_2021: actin cytoskeleton reorganization (GO:0031532)"
$ws.Range("C10").Value = [double]"0.01926694733395609"
$ws.Range("D10").Value = "ns"
$ws.Range("E10").Value = "ns"
$ws.Range("A11").Value = 293
$ws.Range("B11").Value = "GO_Biological_Process_2021: extracellular matrix organization (GO:0030198)"
$ws.Range("C11").Value = "ns"
$ws.Range("D11").Value = "ns"
$ws.Range("E11").Value = [double]"5.017322667628349e-20"
$ws.Range("A12").Value = 294
$ws.Range("B12").Value = "GO_Biological_Process_2021: external encapsulating structure organization (GO:0045229)"
$ws.Range("C12").Value = "ns"
$ws.Range("D12").Value = "ns"
$ws.Range("E12").Value = [double]"3.119685457737965e-21"
$ws.Range("A13").Value = 296
$ws.Range("B13").Value = "GO_Biological_Process_2021: extracellular structure organization (GO:0043062)"
$ws.Range("C13").Value = "ns"
$ws.Range("D13").Value = "ns"
$ws.Range("E13").Value = [double]"3.119685457737965e-21"
$ws.Range("A14").Value = 428
$ws.Range("B14").Value = "KEGG_2021_Human: ECM-receptor interaction"
$ws.Range("C14").Value = "ns"
$ws.Range("D14").Value = "ns"
$ws.Range("E14").Value = [double]"2.790153639676453e-08"
$ws.Range("A15").Value = 574
$ws.Range("B15").Value = "GO_Cellular_Component_2018: actin cytoskeleton (GO:0015629)"
$ws.Range("C15").Value = "ns"
$ws.Range("D15").Value = [double]"2.577279555364868e-05"
$ws.Range("E15").Value = "ns"
$ws.Range("A16").Value = 826
$ws.Range("B16").Value = "GO_Biological_Process_2021: collagen fibril organization (GO:0030199)"
$ws.Range("C16").Value = "ns"
$ws.Range("D16").Value = "ns"
$ws.Range("E16").Value = [double]"1.099185615104196e-19"
$ws.Range("A17").Value = 938
$ws.Range("B17").Value = "GO_Cellular_Component_2018: stress fiber (GO:0001725)"
$ws.Range("C17").Value = "ns"
$ws.Range("D17").Value = [double]"0.0045252846192715"
$ws.Range("E17").Value = "ns"
